$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers for data rows (2-51)
$rowNums = @(
    2,
    3,
    4,
    5,
    6,
    7,
    8,
    9,
    10,
    11,
    12,
    13,
    14,
    15,
    16,
    17,
    18,
    19,
    20,
    21,
    22,
    23,
    24,
    25,
    26,
    27,
    28,
    29,
    30,
    31,
    32,
    33,
    34,
    35,
    36,
    37,
    38,
    39,
    40,
    41,
    42,
    43,
    44,
    45,
    46,
    47,
    48,
    49,
    50,
    51
)

# Column B: Coin name
$coinNames = @(
    'Bitcoin',
    'Ethereum',
    'TetherUSD',
    'BNB',
    'Solana',
    'USDC',
    'XRP',
    'LidoStakedEther',
    'Dogecoin',
    'TRON',
    'Toncoin',
    'Cardano',
    'Avalanche',
    'ShibaInu',
    'WrappedliquidstakedEther2.0',
    'WrappedBTC',
    'WrappedEther',
    'Chainlink',
    'Polkadot',
    'BitcoinCash',
    'Uniswap',
    'Dai',
    'Litecoin',
    'SuiNetwork',
    'Aptos',
    'Binance-PegBSC-USD',
    'InternetComputer(DFINITY)',
    'Bittensor',
    'PEPE',
    'Fetch.AI',
    'Kaspa',
    'PancakeSwap',
    'ImmutableX',
    'FirstDigitalUSD',
    'RenderToken',
    'NEARProtocol',
    'Stacks',
    'PolygonEcosystemToken',
    'EthereumClassic',
    'Monero',
    'USDe',
    'OKB',
    'Aave',
    'dogwifhat',
    'Filecoin',
    'Hedera',
    'InjectiveProtocol',
    'Mantle',
    'Stellar',
    'VeChain'
)

# Column C: Link
$coinLinks = @(
    'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc',
    'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth',
    'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt',
    'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb',
    'https://coinranking.com/coin/zNZHO_Sjf+solana-sol',
    'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc',
    'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp',
    'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth',
    'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge',
    'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx',
    'https://coinranking.com/coin/67YlI0K1b+toncoin-ton',
    'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada',
    'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax',
    'https://coinranking.com/coin/xz24e0BjL+shibainu-shib',
    'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth',
    'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc',
    'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',
    'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link',
    'https://coinranking.com/coin/25W7FG7om+polkadot-dot',
    'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',
    'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni',
    'https://coinranking.com/coin/MoTuySvg7+dai-dai',
    'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc',
    'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui',
    'https://coinranking.com/coin/HGYj5JCv5+aptos-apt',
    'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd',
    'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp',
    'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao',
    'https://coinranking.com/coin/03WI8NQPF+pepe-pepe',
    'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet',
    'https://coinranking.com/coin/V8GxkwWow+kaspa-kas',
    'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake',
    'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',
    'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd',
    'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render',
    'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',
    'https://coinranking.com/coin/mMPrMcB7+stacks-stx',
    'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol',
    'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc',
    'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',
    'https://coinranking.com/coin/exbfr2U-0+usde-usde',
    'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb',
    'https://coinranking.com/coin/ixgUfzmLR+aave-aave',
    'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif',
    'https://coinranking.com/coin/ymQub4fuB+filecoin-fil',
    'https://coinranking.com/coin/jad286TjB+hedera-hbar',
    'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj',
    'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt',
    'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm',
    'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
)

# Column D: Price (kept as text, matching source formatting)
$coinPrices = @(
    '61.476.14',
    '2.383.89',
    '1.00',
    '553.59',
    '140.19',
    '1.00',
    '0.525',
    '2.385.44',
    '0.109',
    '0.158',
    '5.37',
    '0.353',
    '25.69',
    '0.0000174',
    '2.813.51',
    '61.437.59',
    '2.383.46',
    '10.98',
    '4.17',
    '321.34',
    '6.70',
    '1.00',
    '64.33',
    '1.75',
    '8.90',
    '0.999',
    '8.22',
    '523.66',
    '0.0₃0911',
    '1.41',
    '0.149',
    '1.85',
    '1.52',
    '0.999',
    '5.57',
    '4.72',
    '1.88',
    '0.379',
    '18.55',
    '145.84',
    '1.00',
    '41.39',
    '148.12',
    '2.19',
    '3.61',
    '0.0526',
    '19.96',
    '0.584',
    '0.0909',
    '0.0226'
)

# Column E: Volume(1h)
$coinVolumes = @(
    '  +1.13%  ',
    '  +1.29%  ',
    '  +0.01%  ',
    '  +2.74%  ',
    '  +2.28%  ',
    '  +0.00%  ',
    '  +1.37%  ',
    '  +1.39%  ',
    '  +4.23%  ',
    '  +2.34%  ',
    '  +2.57%  ',
    '  +4.03%  ',
    '  +3.94%  ',
    '  +7.93%  ',
    '  +1.28%  ',
    '  +1.34%  ',
    '  +1.17%  ',
    '  +4.13%  ',
    '  +3.06%  ',
    '  +2.04%  ',
    '  +2.43%  ',
    '  +0.26%  ',
    '  +2.03%  ',
    '  -7.05%  ',
    '  +6.17%  ',
    '  -0.06%  ',
    '  +3.89%  ',
    '  +4.59%  ',
    '  +2.54%  ',
    '  +1.73%  ',
    '  +3.05%  ',
    '  +4.48%  ',
    '  +0.11%  ',
    '  +0.02%  ',
    '  +6.86%  ',
    '  +3.87%  ',
    '  +6.31%  ',
    '  +2.29%  ',
    '  +0.67%  ',
    '  +5.26%  ',
    '  -0.01%  ',
    '  +3.29%  ',
    '  +7.06%  ',
    '  +4.80%  ',
    '  +3.56%  ',
    '  +3.47%  ',
    '  +2.72%  ',
    '  +3.06%  ',
    '  +1.59%  ',
    '  +2.02%  '
)

# Force column D to Text format so values like "1.00" or "140.19" are not
# auto-converted to numbers by Excel, matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 2).Value = $coinNames[$i]
    $ws.Cells.Item($r, 3).Value = $coinLinks[$i]
    $ws.Cells.Item($r, 4).Value = $coinPrices[$i]
    $ws.Cells.Item($r, 5).Value = $coinVolumes[$i]
}
